$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("r0")

$ws.Range("A35").Value = "ExoT_r0_script_14v0"
$ws.Range("B35").Value = "Issadore Alzheimer's DLB Project"
$ws.Range("C35").Value = "5 mL"
$ws.Range("D35").Value = "5 mL"
$ws.Range("E35").Value = "1 hour"
$ws.Range("F35").Value = "0.5 mL"
$ws.Range("G35").Value = "1 mL/hr"
$ws.Range("H35").Value = "15 mL/hr"
$ws.Range("I35").Value = "700-700-700"
$ws.Range("J35").Value = "0 mins"
$ws.Range("K35").Value = "N"
